$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column F, matching style of existing header row (A1:E1):
# bold font + centered horizontal alignment (same visual style as style index 1).
$ws.Range("F1").Value = "scenario"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").HorizontalAlignment = -4108

# Fill F2:F101 with the scenario label "S5"
$ws.Range("F2:F101").Value = "S5"
